$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append four new bonus rows after the existing data (rows 2-15), continuing
# the running index in column A and reusing the "Bruno" / "Ronaldo"
# technician pair already present in the sheet.

$newRows = @(
    @{ Row = 16; Idx = 14; OS = 69436857; Tech = "Bruno"   },
    @{ Row = 17; Idx = 15; OS = 69436857; Tech = "Ronaldo" },
    @{ Row = 18; Idx = 16; OS = 69445426; Tech = "Bruno"   },
    @{ Row = 19; Idx = 17; OS = 69445426; Tech = "Ronaldo" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Idx          # A: running index
    $ws.Cells.Item($row, 2).Value = 45818            # B: Data (10/06/2025)
    $ws.Cells.Item($row, 3).Value = $r.OS            # C: Ordem de Serviço
    $ws.Cells.Item($row, 4).Value = 5                # D: Bonificação (R$)
    $ws.Cells.Item($row, 5).Value = $r.Tech          # E: Técnico
}

# Format the new date cells as short dates (built-in numFmtId 14).
$ws.Range("B16:B19").NumberFormat = "mm-dd-yy"

# Widen column B to fit the newly added dates.
$ws.Columns("B:B").AutoFit()

# Leave the selection where the last entry was typed.
$ws.Range("E19").Select() | Out-Null
